$wb = $excel.ActiveWorkbook

# --- Update the conversion note on Hoja1!A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$cell = $ws1.Range("A1")
$text = $cell.Value2
$text = $text -replace [regex]::Escape("1000 Bs = 7.35 = 29955.88 pesos"), "1000 Bs = 7.27 = 29587.13 pesos"
$text = $text -replace [regex]::Escape("29955.88 pesos = 7.33 = 960.17 Bs"), "29587.13 pesos = 7.24 = 945.92 Bs"
$cell.Value = $text

# --- Update the rate values on the "tasas" sheet ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 137.5
$ws2.Range("O10").Value = 4068.23
$ws2.Range("N12").Value = 4084.99
$ws2.Range("O12").Value = 130.6
